# Generate Report for Handoff
#
# The localization status report is regenerated: items that were previously
# "Handed back: in sync with en-US" are now "Ready for handoff", the
# handoff timestamps are refreshed, and the (now shorter) status/datetime
# columns are narrowed to fit.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: status + generate-date columns ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-28 13:00:21"

# --- zh-cn sheet: status + handoff datetime ---
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-28 13:00:16"

# --- de-de sheet: status + handoff datetime ---
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-28 13:00:21"

# --- Narrow the status/datetime columns to match the shorter content ---
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33
